# EPBDS-9540 Hide internal methods. Fix collision resolving
#
# The "runMain()" test-steps table on the "Rules" sheet contained an
# internal/diagnostic step (Step3 = "$Step1.toPlain()") that exposed an
# internal method. Remove that whole row so the remaining steps shift up
# and renumber naturally (Step4..Step19 become Step3..Step18 in position,
# keeping their original labels/formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure we're working on the "Rules" sheet, which holds the
# runMain() steps table.
$ws = $wb.Worksheets.Item("Rules")
$ws.Activate()

# Row 49 is: B49 = "Step3", C49 = "= $Step1.toPlain()"
# Select it the way a user would before deleting, then remove the entire
# row, shifting everything below it up by one.
$ws.Range("A49:XFD49").Select()
$ws.Rows.Item(49).Delete()
